# Updates the cryptos list (Price / Volume(1h) columns) with freshly scraped values.
# Values that look numeric are written with a leading apostrophe and the style is
# reset to 'Normal' afterwards so the cell keeps its original plain-text storage
# (matching the source data, which stores prices such as "0.999" or "213.31" as text,
# not as numbers) instead of Excel silently re-typing it as a Number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.403.96' }
    @{ Cell = 'E2'; Value = '  -1.70%  ' }
    @{ Cell = 'D3'; Value = '1.656.10' }
    @{ Cell = 'E3'; Value = '  -0.47%  ' }
    @{ Cell = 'D4'; Value = '0.999' }
    @{ Cell = 'E4'; Value = '  -0.10%  ' }
    @{ Cell = 'D5'; Value = '213.31' }
    @{ Cell = 'E6'; Value = '  -0.35%  ' }
    @{ Cell = 'E7'; Value = '  -0.11%  ' }
    @{ Cell = 'D8'; Value = '23.65' }
    @{ Cell = 'E8'; Value = '  +0.89%  ' }
    @{ Cell = 'D9'; Value = '0.259' }
    @{ Cell = 'E9'; Value = '  -0.62%  ' }
    @{ Cell = 'E10'; Value = '  -1.38%  ' }
    @{ Cell = 'D11'; Value = '0.0875' }
    @{ Cell = 'E11'; Value = '  -0.43%  ' }
    @{ Cell = 'D12'; Value = '1.890.64' }
    @{ Cell = 'E12'; Value = '  -0.44%  ' }
    @{ Cell = 'D13'; Value = '1.661.39' }
    @{ Cell = 'E13'; Value = '  -0.15%  ' }
    @{ Cell = 'E14'; Value = '  -1.62%  ' }
    @{ Cell = 'E15'; Value = '  +3.21%  ' }
    @{ Cell = 'D16'; Value = '65.63' }
    @{ Cell = 'E16'; Value = '  -0.67%  ' }
    @{ Cell = 'D17'; Value = '27.401.90' }
    @{ Cell = 'E17'; Value = '  -1.52%  ' }
    @{ Cell = 'D18'; Value = '231.74' }
    @{ Cell = 'E18'; Value = '  -6.82%  ' }
    @{ Cell = 'D19'; Value = '0.0₃0725' }
    @{ Cell = 'E19'; Value = '  -0.87%  ' }
    @{ Cell = 'E20'; Value = '  -0.05%  ' }
    @{ Cell = 'E21'; Value = '  -0.05%  ' }
    @{ Cell = 'E22'; Value = '  -2.40%  ' }
    @{ Cell = 'D23'; Value = '9.36' }
    @{ Cell = 'E23'; Value = '  +0.23%  ' }
    @{ Cell = 'D24'; Value = '2.03' }
    @{ Cell = 'E24'; Value = '  -1.16%  ' }
    @{ Cell = 'D25'; Value = '147.51' }
    @{ Cell = 'E25'; Value = '  +0.44%  ' }
    @{ Cell = 'E26'; Value = '  -1.28%  ' }
    @{ Cell = 'E27'; Value = '  -2.33%  ' }
    @{ Cell = 'E28'; Value = '  -0.08%  ' }
    @{ Cell = 'D29'; Value = '0.111' }
    @{ Cell = 'E29'; Value = '  -0.39%  ' }
    @{ Cell = 'D30'; Value = '0.0497' }
    @{ Cell = 'E30'; Value = '  -0.61%  ' }
    @{ Cell = 'E31'; Value = '  -4.30%  ' }
    @{ Cell = 'E32'; Value = '  -1.58%  ' }
    @{ Cell = 'D33'; Value = '1.430.25' }
    @{ Cell = 'E33'; Value = '  -0.57%  ' }
    @{ Cell = 'E34'; Value = '  -0.03%  ' }
    @{ Cell = 'E35'; Value = '  +0.45%  ' }
    @{ Cell = 'E36'; Value = '  -0.69%  ' }
    @{ Cell = 'D37'; Value = '0.909' }
    @{ Cell = 'E37'; Value = '  -2.37%  ' }
    @{ Cell = 'D38'; Value = '0.572' }
    @{ Cell = 'E38'; Value = '  -1.62%  ' }
    @{ Cell = 'E39'; Value = '  +0.04%  ' }
    @{ Cell = 'E40'; Value = '  -0.38%  ' }
    @{ Cell = 'D41'; Value = '0.999' }
    @{ Cell = 'E41'; Value = '  -0.10%  ' }
    @{ Cell = 'E42'; Value = '  +1.76%  ' }
    @{ Cell = 'D43'; Value = '64.96' }
    @{ Cell = 'E43'; Value = '  -6.69%  ' }
    @{ Cell = 'E44'; Value = '  -0.21%  ' }
    @{ Cell = 'D45'; Value = '0.794' }
    @{ Cell = 'E45'; Value = '  +0.70%  ' }
    @{ Cell = 'D46'; Value = '1.798.65' }
    @{ Cell = 'E46'; Value = '  -0.35%  ' }
    @{ Cell = 'E47'; Value = '  -1.29%  ' }
    @{ Cell = 'D48'; Value = '87.98' }
    @{ Cell = 'E48'; Value = '  -1.37%  ' }
    @{ Cell = 'E49'; Value = '  -3.09%  ' }
    @{ Cell = 'E50'; Value = '  -0.48%  ' }
    @{ Cell = 'D51'; Value = '7.71' }
    @{ Cell = 'E51'; Value = '  -1.04%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $value = $u.Value
    $looksNumeric = $value -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    if ($looksNumeric) {
        # Force text storage (apostrophe prefix), then strip the resulting
        # quote-prefix style so the cell ends up unstyled, like the original.
        $cell.Value = "'" + $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}
